$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gdf9"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.7495306666666667
$ws.Range("H2").Value = 2.248592
$ws.Range("I2").Value = 0.08809705645632541
$ws.Range("J2").Value = 0.08809705645632541
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.510190333333334
$ws.Range("N2").Value = 10.530571
$ws.Range("O2").Value = 0.8478537661184122
$ws.Range("P2").Value = 0.8478537661184122
$ws.Range("Q2").Value = 2.630995300670222
$ws.Range("R2").Value = 23.678957706032
$ws.Range("S2").Value = 0.07469342110044189
$ws.Range("T2").Value = 0.07469342110044189

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gdf9"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.7495306666666667
$ws.Range("H3").Value = 2.248592
$ws.Range("I3").Value = 0.08809705645632541
$ws.Range("J3").Value = 0.08809705645632541
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.629899
$ws.Range("N3").Value = 1.889697
$ws.Range("O3").Value = 0.1521462338815877
$ws.Range("P3").Value = 0.1521462338815877
$ws.Range("Q3").Value = 0.4721286174026667
$ws.Range("R3").Value = 4.249157556624
$ws.Range("S3").Value = 0.01340363535588352
$ws.Range("T3").Value = 0.01340363535588353

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdf9"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.763321333333333
$ws.Range("H4").Value = 17.289964
$ws.Range("I4").Value = 0.6773994280135454
$ws.Range("J4").Value = 0.6773994280135454
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.510190333333334
$ws.Range("N4").Value = 10.530571
$ws.Range("O4").Value = 0.8478537661184122
$ws.Range("P4").Value = 0.8478537661184122
$ws.Range("Q4").Value = 20.23035483216044
$ws.Range("R4").Value = 182.073193489444
$ws.Range("S4").Value = 0.5743356562077427
$ws.Range("T4").Value = 0.5743356562077427

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdf9"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.763321333333333
$ws.Range("H5").Value = 17.289964
$ws.Range("I5").Value = 0.6773994280135454
$ws.Range("J5").Value = 0.6773994280135454
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.629899
$ws.Range("N5").Value = 1.889697
$ws.Range("O5").Value = 0.1521462338815877
$ws.Range("P5").Value = 0.1521462338815877
$ws.Range("Q5").Value = 3.630310344545333
$ws.Range("R5").Value = 32.672793100908
$ws.Range("S5").Value = 0.1030637718058026
$ws.Range("T5").Value = 0.1030637718058026

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Gdf9"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1060023333333333
$ws.Range("H6").Value = 0.318007
$ws.Range("I6").Value = 0.01245912136684053
$ws.Range("J6").Value = 0.01245912136684053
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.510190333333334
$ws.Range("N6").Value = 10.530571
$ws.Range("O6").Value = 0.8478537661184122
$ws.Range("P6").Value = 0.8478537661184122
$ws.Range("Q6").Value = 0.3720883657774444
$ws.Range("R6").Value = 3.348795291997
$ws.Range("S6").Value = 0.01056351297340212
$ws.Range("T6").Value = 0.01056351297340212

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Gdf9"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1060023333333333
$ws.Range("H7").Value = 0.318007
$ws.Range("I7").Value = 0.01245912136684053
$ws.Range("J7").Value = 0.01245912136684053
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.629899
$ws.Range("N7").Value = 1.889697
$ws.Range("O7").Value = 0.1521462338815877
$ws.Range("P7").Value = 0.1521462338815877
$ws.Range("Q7").Value = 0.06677076376433333
$ws.Range("R7").Value = 0.600936873879
$ws.Range("S7").Value = 0.001895608393438406
$ws.Range("T7").Value = 0.001895608393438406

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gdf9"
$ws.Range("C8").Value = "Bmpr1b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.889156
$ws.Range("H8").Value = 5.667468
$ws.Range("I8").Value = 0.2220443941632887
$ws.Range("J8").Value = 0.2220443941632887
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.510190333333334
$ws.Range("N8").Value = 10.530571
$ws.Range("O8").Value = 0.8478537661184122
$ws.Range("P8").Value = 0.8478537661184122
$ws.Range("Q8").Value = 6.631297129358667
$ws.Range("R8").Value = 59.68167416422801
$ws.Range("S8").Value = 0.1882611758368255
$ws.Range("T8").Value = 0.1882611758368256

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gdf9"
$ws.Range("C9").Value = "Bmpr1b"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.889156
$ws.Range("H9").Value = 5.667468
$ws.Range("I9").Value = 0.2220443941632887
$ws.Range("J9").Value = 0.2220443941632887
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.629899
$ws.Range("N9").Value = 1.889697
$ws.Range("O9").Value = 0.1521462338815877
$ws.Range("P9").Value = 0.1521462338815877
$ws.Range("Q9").Value = 1.189977475244
$ws.Range("R9").Value = 10.709797277196
$ws.Range("S9").Value = 0.03378321832646317
$ws.Range("T9").Value = 0.03378321832646319
